$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.101290106773376
$ws.Range("B1").Value = 2.282293796539307
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.503005027770996
$ws.Range("E1").Value = 0.9801836013793945
